$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.118.25"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "'1.647.96"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").Value = "'216.28"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "'0.5055"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").Value = "'1.014"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").Value = "'0.2588"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").Value = "'0.06454"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'19.53"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "'0.07763"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "'1.657.05"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'4.264"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "'1.876.49"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "'0.5486"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "'0.0₅7953"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "'63.86"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "'26.151.52"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "'1.014"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "'204.17"
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("D21").Value = "'4.324"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").Value = "'10.04"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").Value = "'5.988"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").Value = "'1.016"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "'1.962"
$ws.Range("E25").Value = "  +11.70%  "
$ws.Range("D26").Value = "'142.39"
$ws.Range("D27").Value = "'0.1161"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'15.75"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").Value = "'6.765"
$ws.Range("E29").Value = "  -2.92%  "
$ws.Range("D30").Value = "'0.05073"
$ws.Range("E30").Value = "  -3.71%  "
$ws.Range("D31").Value = "'1.247"
$ws.Range("D32").Value = "'3.267"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").Value = "'3.207"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").Value = "'1.550"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "'2.353"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").Value = "'0.9019"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("D37").Value = "'2.634"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("D38").Value = "'0.5665"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").Value = "'1.154.22"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").Value = "'0.01578"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").Value = "'2.582"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "'1.014"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").Value = "'5.683"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "'0.8194"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D45").Value = "'100.16"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "'1.786.14"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("E47").Value = "  +4.23%  "
$ws.Range("D48").Value = "'0.4555"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "'1.016"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "'55.15"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").Value = "'0.05049"
$ws.Range("E51").Value = "  -0.93%  "
